$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1:F1 hold strings that look like a date / numbers ("2018-12-02", "12",
# "22.296"); force them to Text up front (as one range op, so it only ever
# allocates a single extra cell style) so Excel's smart-entry parser doesn't
# convert them into a date serial / real numbers.
$ws.Range("D1:F1").NumberFormat = "@"

$ws.Range("A1").Value = "Ford"
$ws.Range("B1").Value = "Isaiah"
$ws.Range("C1").Value = "WR"
$ws.Range("D1").Value = "2018-12-02"
$ws.Range("E1").Value = "12"
$ws.Range("F1").Value = "22.296"
$ws.Range("G1").Value = "MIA"
$ws.Range("H1").Value = ""
$ws.Range("I1").Value = "BUF"
$ws.Range("J1").Value = "W 21-17"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0
